$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap Name/Symbol for rows 24 and 25 (Kaspa <-> Ethereum Classic)
$ws.Range("A24").Value = "Ethereum Classic"
$ws.Range("B24").Value = "ETC-USD"
$ws.Range("A25").Value = "Kaspa"
$ws.Range("B25").Value = "KAS-USD"

# Update Market Cap (column C) values for rows 2-26 with refreshed figures
$ws.Range("C2").Value = 722372549397.9398
$ws.Range("C3").Value = 245521447938.1596
$ws.Range("C4").Value = 35021155768.29485
$ws.Range("C5").Value = 33120120499.31993
$ws.Range("C6").Value = 23975012550.49981
$ws.Range("C7").Value = 13351213137.83818
$ws.Range("C8").Value = 10669931425.47692
$ws.Range("C9").Value = 8951670740.955879
$ws.Range("C10").Value = 8152034859.065285
$ws.Range("C11").Value = 7896231715.013188
$ws.Range("C12").Value = 7334370645.503334
$ws.Range("C13").Value = 7108959453.460295
$ws.Range("C14").Value = 6375380161.262115
$ws.Range("C15").Value = 6029360761.003362
$ws.Range("C16").Value = 5072080315.086798
$ws.Range("C17").Value = 4729853488.20554
$ws.Range("C18").Value = 4362049864.072336
$ws.Range("C19").Value = 3790164544.140406
$ws.Range("C20").Value = 3452298824.828917
$ws.Range("C21").Value = 3271220148.296094
$ws.Range("C22").Value = 3260350550.203758
$ws.Range("C23").Value = 3044044602.572844
$ws.Range("C24").Value = 2703391028.716319
$ws.Range("C25").Value = 2690569435.013835
$ws.Range("C26").Value = 2314128781.164953
